$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("R4").Value = 2023
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4104)
